$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update AttackRate (column H) values for each enemy row
$ws.Range("H2").Value = 10
$ws.Range("H3").Value = 15
$ws.Range("H4").Value = 200

# Update the active cell selection to K7
$ws.Range("K7").Select()
